$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.042550952377426
$ws.Cells.Item(2, 4).Value = 1.04626258578844
$ws.Cells.Item(2, 5).Value = 1.058171073134544
$ws.Cells.Item(2, 6).Value = 1.065436855929432
$ws.Cells.Item(2, 9).Value = 1.045247703823105
$ws.Cells.Item(2, 10).Value = 1.047626167743156
$ws.Cells.Item(2, 11).Value = 1.049028513544368
$ws.Cells.Item(2, 12).Value = 1.060904037257699
$ws.Cells.Item(2, 13).Value = 1.068150099556459
$ws.Cells.Item(2, 14).Value = 1.019715539324603
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.043570650321521
$ws.Cells.Item(3, 4).Value = 1.047063344749364
$ws.Cells.Item(3, 5).Value = 1.059318491080247
$ws.Cells.Item(3, 6).Value = 1.066625080528694
$ws.Cells.Item(3, 9).Value = 1.045584376023714
$ws.Cells.Item(3, 10).Value = 1.048292001340359
$ws.Cells.Item(3, 11).Value = 1.049640907810585
$ws.Cells.Item(3, 12).Value = 1.061864590060655
$ws.Cells.Item(3, 13).Value = 1.06915279106903
$ws.Cells.Item(3, 14).Value = 1.019941468569128
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.044230277128939
$ws.Cells.Item(4, 4).Value = 1.047581258023352
$ws.Cells.Item(4, 5).Value = 1.060061707928994
$ws.Cells.Item(4, 6).Value = 1.067394587894369
$ws.Cells.Item(4, 9).Value = 1.045800825015711
$ws.Cells.Item(4, 10).Value = 1.048722059370135
$ws.Cells.Item(4, 11).Value = 1.050036283447408
$ws.Cells.Item(4, 12).Value = 1.062486308186568
$ws.Cells.Item(4, 13).Value = 1.069801668228632
$ws.Cells.Item(4, 14).Value = 1.020087285160783
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.044507539970328
$ws.Cells.Item(5, 4).Value = 1.047798932998742
$ws.Cells.Item(5, 5).Value = 1.060374338818515
$ws.Cells.Item(5, 6).Value = 1.067718244028635
$ws.Cells.Item(5, 9).Value = 1.045891484848253
$ws.Cells.Item(5, 10).Value = 1.048902668774837
$ws.Cells.Item(5, 11).Value = 1.050202287483889
$ws.Cells.Item(5, 12).Value = 1.062747721028903
$ws.Cells.Item(5, 13).Value = 1.070074472973225
$ws.Cells.Item(5, 14).Value = 1.020148496641947
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.044554091032739
$ws.Cells.Item(6, 4).Value = 1.047835478303568
$ws.Cells.Item(6, 5).Value = 1.060426841635982
$ws.Cells.Item(6, 6).Value = 1.067772596433889
$ws.Cells.Item(6, 9).Value = 1.045906687367878
$ws.Cells.Item(6, 10).Value = 1.048932982921667
$ws.Cells.Item(6, 11).Value = 1.050230147862786
$ws.Cells.Item(6, 12).Value = 1.062791615900432
$ws.Cells.Item(6, 13).Value = 1.070120279084569
$ws.Cells.Item(6, 14).Value = 1.020158769049306
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.044233982102223
$ws.Cells.Item(7, 4).Value = 1.047584166826633
$ws.Cells.Item(7, 5).Value = 1.060065884600513
$ws.Cells.Item(7, 6).Value = 1.067398911992935
$ws.Cells.Item(7, 9).Value = 1.045802037734493
$ws.Cells.Item(7, 10).Value = 1.048724473416661
$ws.Cells.Item(7, 11).Value = 1.05003850243421
$ws.Cells.Item(7, 12).Value = 1.062489801030151
$ws.Cells.Item(7, 13).Value = 1.069805313392589
$ws.Cells.Item(7, 14).Value = 1.020088103424623
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.042895602722743
$ws.Cells.Item(8, 4).Value = 1.046533253692985
$ws.Cells.Item(8, 5).Value = 1.058558691136786
$ws.Cells.Item(8, 6).Value = 1.065838288378842
$ws.Cells.Item(8, 9).Value = 1.045361773704768
$ws.Cells.Item(8, 10).Value = 1.047851350966141
$ws.Cells.Item(8, 11).Value = 1.049235658083749
$ws.Cells.Item(8, 12).Value = 1.06122862429553
$ws.Cells.Item(8, 13).Value = 1.068488949825655
$ws.Cells.Item(8, 14).Value = 1.019791970701686
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.040535771523545
$ws.Cells.Item(9, 4).Value = 1.044679655422271
$ws.Cells.Item(9, 5).Value = 1.055908631212894
$ws.Cells.Item(9, 6).Value = 1.063093202255995
$ws.Cells.Item(9, 9).Value = 1.044575248794527
$ws.Cells.Item(9, 10).Value = 1.046306819534248
$ws.Cells.Item(9, 11).Value = 1.047814176609168
$ws.Cells.Item(9, 12).Value = 1.059007603492069
$ws.Cells.Item(9, 13).Value = 1.06616985972583
$ws.Cells.Item(9, 14).Value = 1.019267281290216
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.03896156617775
$ws.Cells.Item(10, 4).Value = 1.04344276318319
$ws.Cells.Item(10, 5).Value = 1.054145800604722
$ws.Cells.Item(10, 6).Value = 1.061266430027801
$ws.Cells.Item(10, 9).Value = 1.044043690421357
$ws.Cells.Item(10, 10).Value = 1.045273112966081
$ws.Cells.Item(10, 11).Value = 1.046861978717814
$ws.Cells.Item(10, 12).Value = 1.05752779498433
$ws.Cells.Item(10, 13).Value = 1.064624118769315
$ws.Cells.Item(10, 14).Value = 1.018915566648942
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.038279679125258
$ws.Cells.Item(11, 4).Value = 1.042906903466231
$ws.Cells.Item(11, 5).Value = 1.053383387625927
$ws.Cells.Item(11, 6).Value = 1.060476188933437
$ws.Cells.Item(11, 9).Value = 1.043811809796197
$ws.Cells.Item(11, 10).Value = 1.044824552650367
$ws.Cells.Item(11, 11).Value = 1.046448589429116
$ws.Cells.Item(11, 12).Value = 1.056887222151728
$ws.Cells.Item(11, 13).Value = 1.063954865338497
$ws.Cells.Item(11, 14).Value = 1.018762815331355
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.038026358425932
$ws.Cells.Item(12, 4).Value = 1.04270781984202
$ws.Cells.Item(12, 5).Value = 1.053100328645655
$ws.Cells.Item(12, 6).Value = 1.060182772117219
$ws.Cells.Item(12, 9).Value = 1.043725421649273
$ws.Cells.Item(12, 10).Value = 1.044657793203385
$ws.Cells.Item(12, 11).Value = 1.046294875735446
$ws.Cells.Item(12, 12).Value = 1.05664931339966
$ws.Cells.Item(12, 13).Value = 1.063706283585087
$ws.Cells.Item(12, 14).Value = 1.018706008128319
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.038080698281306
$ws.Cells.Item(13, 4).Value = 1.042750525828754
$ws.Cells.Item(13, 5).Value = 1.053161039653581
$ws.Cells.Item(13, 6).Value = 1.060245705877928
$ws.Cells.Item(13, 9).Value = 1.043743963847616
$ws.Cells.Item(13, 10).Value = 1.044693570197612
$ws.Cells.Item(13, 11).Value = 1.04632785519907
$ws.Cells.Item(13, 12).Value = 1.056700344354084
$ws.Cells.Item(13, 13).Value = 1.063759604828384
$ws.Cells.Item(13, 14).Value = 1.018718196568321
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.038258740325222
$ws.Cells.Item(14, 4).Value = 1.042890447996602
$ws.Cells.Item(14, 5).Value = 1.053359987129495
$ws.Cells.Item(14, 6).Value = 1.06045193268856
$ws.Cells.Item(14, 9).Value = 1.04380467417562
$ws.Cells.Item(14, 10).Value = 1.044810771204496
$ws.Cells.Item(14, 11).Value = 1.046435886722604
$ws.Cells.Item(14, 12).Value = 1.056867555961905
$ws.Cells.Item(14, 13).Value = 1.063934317317993
$ws.Cells.Item(14, 14).Value = 1.01875812102871
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.038368432920299
$ws.Cells.Item(15, 4).Value = 1.042976653172963
$ws.Cells.Item(15, 5).Value = 1.053482583122385
$ws.Cells.Item(15, 6).Value = 1.060579010904339
$ws.Cells.Item(15, 9).Value = 1.043842045707001
$ws.Cells.Item(15, 10).Value = 1.044882963503699
$ws.Cells.Item(15, 11).Value = 1.046502426974654
$ws.Cells.Item(15, 12).Value = 1.056970584312506
$ws.Cells.Item(15, 13).Value = 1.064041964592118
$ws.Cells.Item(15, 14).Value = 1.018782710722138
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.039006815549706
$ws.Cells.Item(16, 4).Value = 1.043478320602808
$ws.Cells.Item(16, 5).Value = 1.054196418434744
$ws.Cells.Item(16, 6).Value = 1.061318891743375
$ws.Cells.Item(16, 9).Value = 1.044059043489296
$ws.Cells.Item(16, 10).Value = 1.04530286222755
$ws.Cells.Item(16, 11).Value = 1.046889391217077
$ws.Cells.Item(16, 12).Value = 1.057570311727463
$ws.Cells.Item(16, 13).Value = 1.064668536180007
$ws.Cells.Item(16, 14).Value = 1.01892569463426
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.039407189913827
$ws.Cells.Item(17, 4).Value = 1.043792929239026
$ws.Cells.Item(17, 5).Value = 1.054644430069827
$ws.Cells.Item(17, 6).Value = 1.061783203268147
$ws.Cells.Item(17, 9).Value = 1.044194701953233
$ws.Cells.Item(17, 10).Value = 1.045565996858194
$ws.Cells.Item(17, 11).Value = 1.047131834128655
$ws.Cells.Item(17, 12).Value = 1.057946556388967
$ws.Cells.Item(17, 13).Value = 1.065061584294421
$ws.Cells.Item(17, 14).Value = 1.019015262452875
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.039640697833324
$ws.Cells.Item(18, 4).Value = 1.043976408374486
$ws.Cells.Item(18, 5).Value = 1.054905834977195
$ws.Cells.Item(18, 6).Value = 1.062054102027373
$ws.Cells.Item(18, 9).Value = 1.044273664030507
$ws.Cells.Item(18, 10).Value = 1.045719386361099
$ws.Cells.Item(18, 11).Value = 1.047273142649068
$ws.Cells.Item(18, 12).Value = 1.058166032296286
$ws.Cells.Item(18, 13).Value = 1.0652908486926
$ws.Cells.Item(18, 14).Value = 1.019067461792984
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.039720314001296
$ws.Cells.Item(19, 4).Value = 1.044038965473719
$ws.Cells.Item(19, 5).Value = 1.054994982124767
$ws.Cells.Item(19, 6).Value = 1.062146484028298
$ws.Cells.Item(19, 9).Value = 1.044300560038876
$ws.Cells.Item(19, 10).Value = 1.045771672537162
$ws.Cells.Item(19, 11).Value = 1.04732130750843
$ws.Cells.Item(19, 12).Value = 1.058240871110437
$ws.Cells.Item(19, 13).Value = 1.065369022978808
$ws.Cells.Item(19, 14).Value = 1.019085252945915
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.039364235986235
$ws.Cells.Item(20, 4).Value = 1.043759177476205
$ws.Cells.Item(20, 5).Value = 1.054596353636451
$ws.Cells.Item(20, 6).Value = 1.061733379417762
$ws.Cells.Item(20, 9).Value = 1.044180164172955
$ws.Cells.Item(20, 10).Value = 1.045537774568823
$ws.Cells.Item(20, 11).Value = 1.047105833090106
$ws.Cells.Item(20, 12).Value = 1.057906186961318
$ws.Cells.Item(20, 13).Value = 1.065019413326781
$ws.Cells.Item(20, 14).Value = 1.019005657229969
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.038206312420237
$ws.Cells.Item(21, 4).Value = 1.042849245548858
$ws.Cells.Item(21, 5).Value = 1.053301398331497
$ws.Cells.Item(21, 6).Value = 1.060391200894551
$ws.Cells.Item(21, 9).Value = 1.043786803602322
$ws.Cells.Item(21, 10).Value = 1.044776262410329
$ws.Cells.Item(21, 11).Value = 1.046404078617896
$ws.Cells.Item(21, 12).Value = 1.056818315528801
$ws.Cells.Item(21, 13).Value = 1.063882868607416
$ws.Cells.Item(21, 14).Value = 1.018746366162635
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.037478062748437
$ws.Cells.Item(22, 4).Value = 1.042276895408847
$ws.Cells.Item(22, 5).Value = 1.052487990098947
$ws.Cells.Item(22, 6).Value = 1.059547978664739
$ws.Cells.Item(22, 9).Value = 1.043537993155394
$ws.Cells.Item(22, 10).Value = 1.044296635267427
$ws.Cells.Item(22, 11).Value = 1.04596191749546
$ws.Cells.Item(22, 12).Value = 1.056134492277478
$ws.Cells.Item(22, 13).Value = 1.063168328256604
$ws.Cells.Item(22, 14).Value = 1.018582942911022
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.037864142452343
$ws.Cells.Item(23, 4).Value = 1.04258033171956
$ws.Cells.Item(23, 5).Value = 1.052919119348038
$ws.Cells.Item(23, 6).Value = 1.059994924544644
$ws.Cells.Item(23, 9).Value = 1.043670033457821
$ws.Cells.Item(23, 10).Value = 1.044550973819953
$ws.Cells.Item(23, 11).Value = 1.046196404612655
$ws.Cells.Item(23, 12).Value = 1.056496984519449
$ws.Cells.Item(23, 13).Value = 1.063547115048362
$ws.Cells.Item(23, 14).Value = 1.018669614258193
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.039383645090084
$ws.Cells.Item(24, 4).Value = 1.043774428527586
$ws.Cells.Item(24, 5).Value = 1.054618077037255
$ws.Cells.Item(24, 6).Value = 1.061755892443092
$ws.Cells.Item(24, 9).Value = 1.044186733680463
$ws.Cells.Item(24, 10).Value = 1.045550527292164
$ws.Cells.Item(24, 11).Value = 1.047117582162461
$ws.Cells.Item(24, 12).Value = 1.057924428109777
$ws.Cells.Item(24, 13).Value = 1.06503846855278
$ws.Cells.Item(24, 14).Value = 1.01900999755283
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.041146016575133
$ws.Cells.Item(25, 4).Value = 1.045159060921626
$ws.Cells.Item(25, 5).Value = 1.056593049964786
$ws.Cells.Item(25, 6).Value = 1.063802291529221
$ws.Cells.Item(25, 9).Value = 1.044779854974237
$ws.Cells.Item(25, 10).Value = 1.046706826322658
$ws.Cells.Item(25, 11).Value = 1.048182464931822
$ws.Cells.Item(25, 12).Value = 1.059581635107459
$ws.Cells.Item(25, 13).Value = 1.06676934209084
$ws.Cells.Item(25, 14).Value = 1.019403265160395